$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "320018398400"
$ws.Range("P2").NumberFormat = "General"

$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "320018398411"
$ws.Range("P3").NumberFormat = "General"

$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "320018398444"
$ws.Range("P4").NumberFormat = "General"
